$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet from "Scanner" to "Session"
$ws.Name = "Session"

# Delete row 2 (the data row), shifting dimension back to A1:F1
$ws.Rows.Item(2).Delete()
